# "Updated instructions for next year"
# - Duplicate the Sprint 2 rubric ("Sheet1 (3)") into a fresh "Sheet1 (4)"
#   tab for next year: reset the "Actual" scores back to match "Possible"
#   (a clean, ungraded template) and drop the old grading-issue notes.
# - Nudge the saved scroll position / selection on the two existing
#   Sprint tabs so the workbook reopens where the instructor left off.

$wb = $excel.ActiveWorkbook

# --- Sheet1 (2): just scrolled down a bit, cell E11 now selected ---
$ws2 = $wb.Worksheets.Item("Sheet1 (2)")
$ws2.Activate()
$ws2.Range("E11").Select()

# --- Sheet1 (3): scrolled further down, no longer the active tab ---
$ws3 = $wb.Worksheets.Item("Sheet1 (3)")
$ws3.Activate()
$ws3.Range("C6").Select()

# --- New Sheet1 (4): copy of Sheet1 (3), reset as a clean template ---
$ws3.Copy($null, $ws3)
$ws4 = $wb.Worksheets.Item($ws3.Index + 1)
$ws4.Name = "Sheet1 (4)"

# Remove the old "Issues" notes from the bottom of the sheet
$ws4.Range("A18:A23").ClearContents()

# Clear the per-row callout notes in column D
$ws4.Range("D6").ClearContents()
$ws4.Range("D10").ClearContents()
$ws4.Range("D13").ClearContents()

# Reset "Actual" (column C) back to match "Possible" (column B) so next
# year's copy starts fresh/ungraded
$ws4.Range("C6").Value = 40
$ws4.Range("C10").Value = 5
$ws4.Range("C13").Value = 5

$ws4.Activate()
$ws4.Range("A18").Select()
